$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.512.71'
$ws.Range("E2").Value = '  -0.05%  '
$ws.Range("D3").Value = '2.503.05'
$ws.Range("E3").Value = '  -0.29%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '574.95'
$ws.Range("E5").Value = '  -0.73%  '
$ws.Range("D6").Value = '166.59'
$ws.Range("E6").Value = '  -0.30%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '0.514'
$ws.Range("E8").Value = '  -0.75%  '
$ws.Range("D9").Value = '2.502.09'
$ws.Range("E9").Value = '  -0.27%  '
$ws.Range("D10").Value = '0.160'
$ws.Range("E10").Value = '  +0.57%  '
$ws.Range("D11").Value = '0.167'
$ws.Range("E11").Value = '  -0.57%  '
$ws.Range("D12").Value = '0.358'
$ws.Range("E12").Value = '  +5.04%  '
$ws.Range("D13").Value = '4.91'
$ws.Range("E13").Value = '  +1.09%  '
$ws.Range("D14").Value = '2.962.68'
$ws.Range("E14").Value = '  -0.06%  '
$ws.Range("D15").Value = '69.308.56'
$ws.Range("E15").Value = '  -0.11%  '
$ws.Range("E16").Value = '  +0.53%  '
$ws.Range("D17").Value = '24.83'
$ws.Range("E17").Value = '  -0.49%  '
$ws.Range("D18").Value = '2.507.33'
$ws.Range("E18").Value = '  +0.13%  '
$ws.Range("D19").Value = '11.23'
$ws.Range("E19").Value = '  -2.28%  '
$ws.Range("D20").Value = '7.55'
$ws.Range("E20").Value = '  -3.28%  '
$ws.Range("D21").Value = '349.94'
$ws.Range("E21").Value = '  -0.46%  '
$ws.Range("D22").Value = '3.91'
$ws.Range("E22").Value = '  -1.45%  '
$ws.Range("E23").Value = '  -1.26%  '
$ws.Range("D25").Value = '70.40'
$ws.Range("E25").Value = '  +1.65%  '
$ws.Range("E26").Value = '  -2.16%  '
$ws.Range("E27").Value = '  -2.57%  '
$ws.Range("D28").Value = '2.636.90'
$ws.Range("E28").Value = '  -0.02%  '
$ws.Range("E29").Value = '  -0.70%  '
$ws.Range("D30").Value = '0.0₃0888'
$ws.Range("E30").Value = '  -2.00%  '
$ws.Range("D31").Value = '7.86'
$ws.Range("E31").Value = '  -0.37%  '
$ws.Range("D32").Value = '460.63'
$ws.Range("E32").Value = '  -3.90%  '
$ws.Range("E33").Value = '  -5.43%  '
$ws.Range("D34").Value = '1.73'
$ws.Range("E34").Value = '  -1.29%  '
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("D36").Value = '158.97'
$ws.Range("E36").Value = '  +4.14%  '
$ws.Range("D37").Value = '0.117'
$ws.Range("E37").Value = '  +0.88%  '
$ws.Range("E38").Value = '  +0.97%  '
$ws.Range("E39").Value = '  -0.81%  '
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("D41").Value = '0.318'
$ws.Range("E41").Value = '  -0.48%  '
$ws.Range("D42").Value = '4.67'
$ws.Range("E42").Value = '  -2.14%  '
$ws.Range("E43").Value = '  -1.94%  '
$ws.Range("D44").Value = '38.20'
$ws.Range("E44").Value = '  +0.12%  '
$ws.Range("D45").Value = '2.21'
$ws.Range("D46").Value = '1.09'
$ws.Range("E46").Value = '  -7.12%  '
$ws.Range("D47").Value = '142.25'
$ws.Range("E47").Value = '  -1.44%  '
$ws.Range("D48").Value = '3.46'
$ws.Range("E48").Value = '  -2.48%  '
$ws.Range("D49").Value = '0.519'
$ws.Range("E49").Value = '  -2.46%  '
$ws.Range("E50").Value = '  +0.37%  '
$ws.Range("D51").Value = '5.77'
$ws.Range("E51").Value = '  +2.71%  '
